$d = $word.ActiveDocument

# 1) Merge the two bold runs "LangGraph " + "Sequential flow" into a single run.
$d.Content.Find.Execute("LangGraph Sequential flow", $true, $false, $false, $false, $false, $true, 1, $false, "LangGraph Sequential flow", 2) | Out-Null

# 2) Append new paragraphs at the end of the document: four blank paragraphs,
#    an "Output:" paragraph, and a paragraph holding the captured run output.
for ($i = 0; $i -lt 5; $i++) {
    $end = $d.Content.End
    $r = $d.Range($end, $end)
    $r.InsertParagraphAfter() | Out-Null
}

$end = $d.Content.End
$r = $d.Range($end, $end)
$r.InsertAfter("Output:") | Out-Null

$end = $d.Content.End
$r = $d.Range($end, $end)
$r.InsertParagraphAfter() | Out-Null

$end = $d.Content.End
$r = $d.Range($end, $end)
$r.InsertAfter('{''topic'': ''Write an essay on greatness of Java'', ''outline'': ''Title: The Greatness of Java: The Time-Tested Programming Language\n\nI. Introduction\n   A. Brief overview of Java\''s history\n      1. Origin in the mid-1990s by Sun Microsystems\n      2. Java\''s growth and evolution\n   B. Importance of programming languages in the tech industry\n   C. Thesis statement: Java\''s greatness stems from its platform independence, robust ecosystem, versatility, and continuous adaptability.\n\nII. Platform Independence\n   A. Explanation of "Write Once, Run Anywhere" (WORA)\n   B. Java Virtual Machine (JVM) as the key enabler\n   C. Portability on diverse hardware and operating systems\n   D. Case studies of applications utilizing Java’s platform-independent nature\n\nIII. Robustness and Security\n   A. Memory management and exception handling\n      1. Automatic garbage collection\n      2. Strong type checking\n   B. Security features\n      1. Bytecode verification\n      2. Security API and access control functionalities\n   C. Real-world examples of Java in secure applications\n      1. Banking systems\n      2. Government and defense\n\nIV. Versatility and Diverse Applications\n   A. Enterprise applications\n      1. Java EE and its role in enterprise architecture\n      2. Examples of large-scale enterprise systems\n   B. Mobile applications\n      1. Java\''s role in Android development\n      2. Case study of popular Android applications developed in Java\n   C. Web development\n      1. Java-based frameworks like Spring\n      2. Use in scalable web services and applications\n   D. Emerging technologies\n      1. Adoption in Internet of Things (IoT)\n      2. Role in cloud computing and big data\n\nV. Extensive Ecosystem and Community Support\n   A. Wide range of libraries and frameworks\n      1. Explanation of open-source contributions\n      2. Popular Java frameworks and tools: Spring, Hibernate, Maven\n   B. Strong and active community\n      1. Java Community Process (JCP)\n      2. Support networks and resources, including forums and educational sites\n   C. Contribution to innovations and industry standards\n\nVI. Continuous Adaptability and Evolution\n   A. Regular updates and enhancements\n      1. Overview of Java\''s major version releases\n      2. Introduction of new features: lambda expressions, streams, modules\n   B. Maintaining relevance in modern software development\n   C. Future prospects of Java in technology trends\n      1. Anticipated developments in AI and machine learning\n      2. Role in future-proof software solutions\n\nVII. Criticisms and Challenges\n   A. Performance considerations against lower-level languages\n   B. Complexity in comparison to newer languages\n   C. Discussion on overcoming these challenges\n\nVIII. Conclusion\n   A. Recap of Java’s core strengths and contributions to technology\n   B. Personal reflection on Java’s impact on the author or industry professionals\n   C. Closing thoughts on the enduring greatness of Java in a rapidly evolving technological landscape\n\nIX. References\n   A. Citing books, articles, and online resources on Java\n   B. Inclusion of quotes and insights from industry experts\n   C. Reference to Java documentation and technical specifications\n\nThis outline provides a comprehensive structure for exploring Java\''s greatness, considering its historical context, technical advantages, and influence on diverse technological sectors.'', ''story'': ''**The Greatness of Java: The Time-Tested Programming Language**\n\nJava, originating in the mid-1990s under the auspices of Sun Microsystems, has transformed from a fledgling language to a cornerstone of modern software development. As tech industries burgeon with numerous programming languages, Java’s enduring relevance stems from its platform independence, robust ecosystem, versatility, and continuous adaptability.\n\nJava’s mantra, “Write Once, Run Anywhere” (WORA), is enabled by the Java Virtual Machine (JVM), which allows portability across diverse hardware and operating systems. This ability has empowered countless applications across sectors, demonstrating Java’s unparalleled flexibility and resilience in evolving tech landscapes.\n\nRobustness and security are hallmarks of Java’s design. Its automatic memory management and stringent type checking minimize common programming errors, while features like bytecode verification and comprehensive security APIs ensure secure application execution. Java’s reliability makes it a preferred choice for critical systems in banking, government, and defense, where security cannot be compromised.\n\nJava’s versatility is evident in its widespread applications. From powering large-scale enterprise systems through Java EE to driving innovation in mobile development, particularly in Android, Java’s impact is ubiquitous. In web development, frameworks like Spring allow the creation of scalable services, while its emergence in IoT, cloud computing, and big data underscores its adaptability to new technologies.\n\nThe Java community bolsters its ecosystem, with extensive libraries, tools, and active forums that drive constant innovation and maintain its industry standard status. Regular updates bring features like lambda expressions and modules, ensuring Java remains at the forefront of software development, with anticipated roles in AI and machine learning.\n\nDespite challenges like performance gaps with lower-level languages and complexities compared to newer syntaxes, Java’s adaptability ensures it overcomes such hurdles. Its sustained evolution assures its place in future-proof software solutions.\n\nIn conclusion, Java’s core strengths, its profound industry impact, and its ability to adapt to ever-changing technological demands highlight its enduring greatness. For professionals and enthusiasts alike, Java’s legacy is not just a testament to its past achievements but a beacon guiding future technological advancements in an incessantly evolving landscape.''}') | Out-Null
